$p = $ppt.ActivePresentation
$w = $ppt.ActiveWindow
$v = $w.View
Write-Host "Before: $($v.Type)"
try {
  $v.Type = 5
  Write-Host "After: $($v.Type)"
} catch { Write-Host "ERR: $_" }

$nm = $p.NotesMaster
$t = $nm.Theme
$tcs = $t.ThemeColorScheme
$c = $tcs.Colors(3)
Write-Host "dk2 color now: $($c.RGB)"
$c.RGB = 1111111
